$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colors = @("Red", "Green", "Blue")
$features = @("Mean", "Median", "Variance", "Std")

$col = 1
foreach ($color in $colors) {
    foreach ($feature in $features) {
        $ws.Cells.Item(1, $col).Value = "$($color)_$($feature)"
        $col++
    }
}
